$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove the leading "序号" column entirely; everything shifts one
#    column to the left (B->A, C->B, D->C, E->D).
$ws.Columns.Item(1).Delete()

# 2) Remove the now-empty/orphan row (old row 7, which only held a stray
#    "福" label with no other data); everything below shifts up by one row.
$ws.Rows.Item(7).Delete()

# 3) Correct the label text for the fungus row (was "酵母样真菌",
#    now just "真菌") - this is the new row 8 after the shifts above.
$ws.Cells.Item(8, 1).Value = "真菌"
